$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Total" (column B) values
$ws.Range("B2").Value = 116037.66
$ws.Range("B3").Value = 956806.1899999999
$ws.Range("B4").Value = 1772410.88
$ws.Range("B5").Value = 2846150.66
$ws.Range("B6").Value = 4457482.68
$ws.Range("B7").Value = 1989051.38

# Update "Evolução Total (%)" (column C) values
$ws.Range("C3").Value = 724.5652230491376
$ws.Range("C4").Value = 85.24241361774634
$ws.Range("C5").Value = 60.58074863544059
$ws.Range("C6").Value = 56.61443164783131
$ws.Range("C7").Value = -55.3772493850722

# Update "Qtd Produtos" (column D) value for 2025
$ws.Range("D7").Value = 1746
